$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enter the new literal values into the pairing matrix data-entry cells.
$ws.Range("D14").Value = 60
$ws.Range("E14").Value = 60
$ws.Range("H15").Value = 360

# B17 previously held the formula "=E14"; the user retyped it as a plain
# literal value (60), so overwrite the formula with a hard-coded value.
$ws.Range("B17").Value = 60

# Update the active selection to reflect where the user left the cursor.
$ws.Range("F18").Select()
